# Applies the commit's spreadsheet edit:
#  - General Ledger sheet: rename the "DESCRIPTION" header (C3) to "ITEMS"
#    (a new "ITEMS" shared string is introduced; "BALANCE" stays as-is).
#  - General Journal / General Ledger: column widths are re-tuned from a
#    uniform 3/4-column layout to a 6/7 explicit-width layout (every column
#    gets its own <col> width instead of being grouped).

$wb = $excel.ActiveWorkbook

$wsJournal = $wb.Worksheets.Item("General Journal")
$wsLedger  = $wb.Worksheets.Item("General Ledger")

# --- General Journal column widths ------------------------------------
# target stored widths: A=9.28515625 B=12 C=28.7109375 D=9.28515625 E=12 F=12.85546875
# (ColumnWidth in this runtime is stored as ColumnWidth + 5/6, so subtract
#  5/6 = 0.8333333333333334 from every target to get the value to assign)
$wsJournal.Columns.Item(1).ColumnWidth = 8.451822916666666
$wsJournal.Columns.Item(2).ColumnWidth = 11.166666666666666
$wsJournal.Columns.Item(3).ColumnWidth = 27.877604166666668
$wsJournal.Columns.Item(4).ColumnWidth = 8.451822916666666
$wsJournal.Columns.Item(5).ColumnWidth = 11.166666666666666
$wsJournal.Columns.Item(6).ColumnWidth = 12.022135416666666

# --- General Ledger column widths --------------------------------------
# target stored widths: A=7.5703125 B=11.5703125 C=21.7109375 D=7.5703125
#                        E=11.5703125 F=12.42578125 G=12.85546875
$wsLedger.Columns.Item(1).ColumnWidth = 6.736979166666667
$wsLedger.Columns.Item(2).ColumnWidth = 10.736979166666666
$wsLedger.Columns.Item(3).ColumnWidth = 20.877604166666668
$wsLedger.Columns.Item(4).ColumnWidth = 6.736979166666667
$wsLedger.Columns.Item(5).ColumnWidth = 10.736979166666666
$wsLedger.Columns.Item(6).ColumnWidth = 11.592447916666666
$wsLedger.Columns.Item(7).ColumnWidth = 12.022135416666666

# --- Header text change --------------------------------------------------
# General Ledger row 3 header: C3 goes from "DESCRIPTION" to "ITEMS"
$wsLedger.Range("C3").Value = "ITEMS"
